$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F ("想去人数") values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 645
$ws1.Range("F6").Value = 9780
$ws1.Range("F7").Value = 882
$ws1.Range("F10").Value = 3358
$ws1.Range("F11").Value = 168
$ws1.Range("F12").Value = 112
$ws1.Range("F13").Value = 34
$ws1.Range("F14").Value = 30
$ws1.Range("F15").Value = 284
$ws1.Range("F16").Value = 531
$ws1.Range("F17").Value = 107
$ws1.Range("F18").Value = 264
$ws1.Range("F19").Value = 1427

# Sheet "全部类型" (all types) - update column F ("想去人数") values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 645
$ws4.Range("F7").Value = 9780
$ws4.Range("F8").Value = 882
$ws4.Range("F11").Value = 3358
$ws4.Range("F12").Value = 168
$ws4.Range("F13").Value = 112
$ws4.Range("F14").Value = 34
$ws4.Range("F15").Value = 30
$ws4.Range("F16").Value = 284
$ws4.Range("F17").Value = 531
$ws4.Range("F18").Value = 107
$ws4.Range("F19").Value = 264
$ws4.Range("F20").Value = 1427
